$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text updates
# M6 is a plain (non-rich-text) string - straightforward whole-value replace.
$ws.Range("M6").Value = "Edward A. Caban"

# A8 ("Volume 30   Number  26") and C9 (the date-range string) are rich-text
# cells built from multiple runs that all happen to share the same font
# (sz 10, Andale WT). Edit only the specific characters that changed (via
# Characters(start,len).Text) instead of overwriting the whole cell, and do
# it right-to-left so earlier offsets aren't invalidated by length changes
# from a prior edit.
$ws.Range("A8").Characters(21, 2).Text = "27"

$ws.Range("C9").Characters(47, 8).Text = "7/9/2023"
$ws.Range("C9").Characters(27, 9).Text = "7/3/2023"

# Table updates (rows 14-30)
$data = @{
    14 = @{ C=4; D=4; E=0; F=14; G=11; H=27.272727272727; I=69; J=75; K=-8; L=-2.816901408450; M=11.290322580645; N=-72.727272727272 }
    15 = @{ C=5; D=9; E=-44.444444444444; F=25; G=33; H=-24.242424242424; I=204; J=210; K=-2.857142857142; L=13.333333333333; M=39.726027397260; N=-45.454545454545 }
    16 = @{ C=101; D=114; E=-11.403508771929; F=423; G=453; H=-6.622516556291; I=2390; J=2505; K=-4.590818363273; L=34.799774393683; M=8.983128134974; N=-71.270585406899 }
    17 = @{ C=184; D=169; E=8.875739644970; F=714; G=664; H=7.530120481927; I=4050; J=3755; K=7.856191744340; L=32.874015748031; M=78.807947019867; N=-13.387510692899 }
    18 = @{ C=40; D=65; E=-38.461538461538; F=192; G=250; H=-23.2; I=1546; J=1552; K=-0.386597938144; L=43.947858472998; M=-6.416464891041; N=-84.030575353785 }
    19 = @{ C=124; D=141; E=-12.056737588652; F=565; G=633; H=-10.742496050552; I=3958; J=4083; K=-3.061474406073; L=24.504561182761; M=72.012168622338; N=4.875463698993 }
    20 = @{ C=74; D=83; E=-10.843373493975; F=374; G=300; H=24.666666666666; I=2743; J=2079; K=31.938431938431; L=106.706857573474; M=162.739463601533; N=-65.374905326937 }
    21 = @{ C=532; D=585; E=-9.059829059829; F=2307; G=2344; H=-1.578498293515; I=14960; J=14259; K=4.916193281436; L=40.443109275253; M=54.817344509986; N=-57.255921597759 }
    22 = @{ C=4; D=3; E=33.333333333333; F=20; G=28; H=-28.571428571428; I=151; J=191; K=-20.942408376963; L=23.770491803278; M=-10.119047619047 }
    23 = @{ C=37; D=31; E=19.354838709677; F=130; G=143; H=-9.090909090909; I=910; J=843; K=7.947805456702; L=50.165016501650; M=69.776119402985 }
    24 = @{ C=269; D=443; E=-39.277652370203; F=1444; G=1606; H=-10.087173100871; I=9217; J=9512; K=-3.101345668629; L=46.767515923566; M=43.768522851349 }
    25 = @{ C=202; D=244; E=-17.213114754098; F=858; G=921; H=-6.840390879478; I=5476; J=5283; K=3.653227332954; L=28.847058823529; M=-5.406806011400 }
    26 = @{ C=8; D=19; E=-57.894736842105; F=40; G=58; H=-31.034482758620; I=337; J=371; K=-9.164420485175; L=11.960132890365 }
    27 = @{ C=18; D=20; E=-10; F=77; G=67; H=14.925373134328; I=550; J=470; K=17.021276595744; L=20.879120879120 }
    28 = @{ C=23; D=16; E=43.75; F=64; G=57; H=12.280701754386; I=210; J=278; K=-24.460431654676; L=-32.475884244373; M=-8.296943231441; N=-68.703427719821 }
    29 = @{ C=17; D=13; E=30.769230769230; F=52; G=48; H=8.333333333333; I=173; J=238; K=-27.310924369747; L=-34.220532319391; M=-9.895833333333; N=-71.452145214521 }
    30 = @{ J=27; K=-55.555555555555 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
